$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updated values
$ws.Range("B2").Value = 0.01514828764759746
$ws.Range("C2").Value = 0.002777888934908601
$ws.Range("D2").Value = 337.1190423067083
$ws.Range("E2").Value = 616238.5361209477
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 616575.673089431

# Row 3 updated values
$ws.Range("B3").Value = 0.6753301551942219
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 0.8054896365839992
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.645393585217082
